$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "67.098.67"
$ws.Range("E2").Value = "  +0.70%  "

Set-TextValue "D3" "2.481.73"
$ws.Range("E3").Value = "  +2.38%  "

$ws.Range("E4").Value = "  +0.14%  "

Set-TextValue "D5" "583.55"
$ws.Range("E5").Value = "  +1.65%  "

Set-TextValue "D6" "170.87"
$ws.Range("E6").Value = "  +3.86%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  +2.40%  "

Set-TextValue "D9" "2.481.04"
$ws.Range("E9").Value = "  +2.51%  "

Set-TextValue "D10" "0.138"
$ws.Range("E10").Value = "  +5.32%  "

$ws.Range("E11").Value = "  +1.73%  "

$ws.Range("E12").Value = "  +3.37%  "

$ws.Range("E13").Value = "  +2.85%  "

Set-TextValue "D14" "25.41"
$ws.Range("E14").Value = "  +2.44%  "

Set-TextValue "D16" "66.983.67"
$ws.Range("E16").Value = "  +1.07%  "

$ws.Range("E17").Value = "  +3.04%  "

Set-TextValue "D18" "2.498.29"
$ws.Range("E18").Value = "  +4.32%  "

Set-TextValue "D19" "11.01"
$ws.Range("E19").Value = "  -0.52%  "

$ws.Range("E20").Value = "  +0.63%  "

Set-TextValue "D21" "348.38"
$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("E22").Value = "  +1.53%  "

$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("E24").Value = "  +2.04%  "

Set-TextValue "D25" "68.46"
$ws.Range("E25").Value = "  +0.22%  "

$ws.Range("E26").Value = "  +3.59%  "

Set-TextValue "D27" "9.30"
$ws.Range("E27").Value = "  +4.41%  "

$ws.Range("E28").Value = "  -0.01%  "

Set-TextValue "D29" "2.621.53"
$ws.Range("E29").Value = "  +4.01%  "

$ws.Range("E30").Value = "  +3.78%  "

Set-TextValue "D31" "511.26"
$ws.Range("E31").Value = "  +4.84%  "

$ws.Range("E32").Value = "  +0.64%  "

$ws.Range("E33").Value = "  +2.87%  "

$ws.Range("E34").Value = "  +0.26%  "

$ws.Range("E35").Value = "  +0.17%  "

Set-TextValue "D36" "160.28"
$ws.Range("E36").Value = "  +2.46%  "

$ws.Range("E37").Value = "  +4.06%  "

$ws.Range("E38").Value = "  +0.92%  "

Set-TextValue "D39" "18.23"
$ws.Range("E39").Value = "  +0.46%  "

$ws.Range("E40").Value = "  +1.20%  "

$ws.Range("E41").Value = "  +2.94%  "

$ws.Range("E42").Value = "  +0.05%  "

$ws.Range("E43").Value = "  +2.42%  "

Set-TextValue "D44" "4.79"
$ws.Range("E44").Value = "  +4.26%  "

$ws.Range("E45").Value = "  +2.86%  "

Set-TextValue "D46" "38.88"
$ws.Range("E46").Value = "  -0.35%  "

Set-TextValue "D47" "142.69"
$ws.Range("E47").Value = "  +3.54%  "

$ws.Range("E48").Value = "  +1.20%  "

Set-TextValue "D49" "0.514"
$ws.Range("E49").Value = "  +1.62%  "

$ws.Range("E50").Value = "  +3.40%  "

$ws.Range("E51").Value = "  +1.65%  "
